# Sample Project / Main.xlsx — "SAVE" edit
#
# Semantic change (per the commit diff): cell B11 on the "Rules" sheet
# changes from the text "R40" to the text "1". Everything else touched by
# the diff (column attribute ordering, styles.xml bookkeeping tags, the
# sharedStrings count/uniqueCount bump) is a side effect of the external
# tool that re-serialized the workbook on that "SAVE" event, not a
# distinct user edit, so we only reproduce the actual cell-content change.
#
# B11 already holds a shared string ("R40"), and the target value "1"
# must stay textual (shared string) rather than become the number 1.
# A plain `Range.Value = "1"` assignment mimics typing into the cell, so
# Excel's normal auto-detection kicks in and stores it as a numeric 1
# (it also happens to keep cell B11's existing style index, but loses the
# text type). To force a literal text "1" while leaving B11's style
# completely untouched, we stage the text in a scratch cell (using a
# formula whose computed result is the text "1"), copy it, and paste only
# the *value* into B11 — copy/paste of a text value preserves its string
# type instead of re-parsing it as user input. The scratch cell (B5) is
# already blank inside the sheet's used range, so clearing it afterwards
# leaves no trace.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("B5")
$target = $ws.Range("B11")

$scratch.Formula = "=TEXT(1,""0"")"
$scratch.Copy()
$target.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$scratch.ClearContents()
$excel.CutCopyMode = $false
